$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab and update the "through" date label
$ws.Name = "Through 2022-07-12"
$ws.Range("I1").Value = "2022 (through 07-12)"

# Update the data that changed for the new date
$ws.Range("I5").Value = 115
$ws.Range("I8").Value = 68
$ws.Range("I14").Value = 874
